# Scheduled-runner refresh: pushes freshly-fetched market price/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) into each crafting
# job's leve sheet. Only numeric market-data columns (H:N) are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 284.33334
$ws.Range("I33").Value = 300.55554
$ws.Range("K33").Value = 300.55554
$ws.Range("M33").Value = -71.55554000000001

$ws.Range("H80").Value = 1341.826
$ws.Range("I80").Value = 2184.818
$ws.Range("J80").Value = 569.0833
$ws.Range("K80").Value = 6554.454000000001
$ws.Range("L80").Value = 1707.2499
$ws.Range("M80").Value = -5556.454000000001
$ws.Range("N80").Value = -3703.2499

$ws.Range("H83").Value = 1341.826
$ws.Range("I83").Value = 2184.818
$ws.Range("J83").Value = 569.0833
$ws.Range("K83").Value = 19663.362
$ws.Range("L83").Value = 5121.7497
$ws.Range("M83").Value = -14671.362
$ws.Range("N83").Value = -15105.7497

$ws.Range("H112").Value = 5435.6577
$ws.Range("J112").Value = 6030.8823
$ws.Range("L112").Value = 18092.6469
$ws.Range("N112").Value = -20308.6469

$ws.Range("H113").Value = 2976.182
$ws.Range("I113").Value = 3047.6
$ws.Range("K113").Value = 3047.6
$ws.Range("M113").Value = 206.4000000000001

$ws.Range("H116").Value = 3718.0908
$ws.Range("I116").Value = 3612.5
$ws.Range("J116").Value = 3999.6667
$ws.Range("K116").Value = 3612.5
$ws.Range("L116").Value = 3999.6667
$ws.Range("M116").Value = -170.5
$ws.Range("N116").Value = -10883.6667

$ws.Range("H125").Value = 1818
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 2090.6667
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 18816.0003
$ws.Range("M125").Value = -6540
$ws.Range("N125").Value = -23736.0003

$ws.Range("H129").Value = 994.5111000000001
$ws.Range("I129").Value = 474.27274
$ws.Range("J129").Value = 1162.8235
$ws.Range("K129").Value = 1422.81822
$ws.Range("L129").Value = 3488.4705
$ws.Range("M129").Value = 3577.18178
$ws.Range("N129").Value = -13488.4705

$ws.Range("H132").Value = 4870.45
$ws.Range("I132").Value = 4600.5264
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 13801.5792
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -11271.5792
$ws.Range("N132").Value = -35057

$ws.Range("H138").Value = 2726.72
$ws.Range("I138").Value = 1825.7778
$ws.Range("J138").Value = 2815.8242
$ws.Range("K138").Value = 5477.3334
$ws.Range("L138").Value = 8447.472600000001
$ws.Range("M138").Value = -337.3334000000004
$ws.Range("N138").Value = -18727.4726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5628260.5
$ws.Range("I32").Value = 6502021
$ws.Range("J32").Value = 21630.584
$ws.Range("K32").Value = 6502021
$ws.Range("L32").Value = 21630.584
$ws.Range("M32").Value = -6501734
$ws.Range("N32").Value = -22204.584

$ws.Range("H61").Value = 10419506
$ws.Range("I61").Value = 12822308
$ws.Range("J61").Value = 7366.6665
$ws.Range("K61").Value = 12822308
$ws.Range("L61").Value = 7366.6665
$ws.Range("M61").Value = -12822096
$ws.Range("N61").Value = -7790.6665

$ws.Range("H74").Value = 1858.0526
$ws.Range("I74").Value = 1266.5
$ws.Range("J74").Value = 3041.158
$ws.Range("K74").Value = 1266.5
$ws.Range("L74").Value = 3041.158
$ws.Range("M74").Value = -392.5
$ws.Range("N74").Value = -4789.157999999999

$ws.Range("H77").Value = 1858.0526
$ws.Range("I77").Value = 1266.5
$ws.Range("J77").Value = 3041.158
$ws.Range("K77").Value = 6332.5
$ws.Range("L77").Value = 15205.79
$ws.Range("M77").Value = -1964.5
$ws.Range("N77").Value = -23941.79

$ws.Range("H122").Value = 112284.664
$ws.Range("J122").Value = 1367.1428
$ws.Range("L122").Value = 4101.428400000001
$ws.Range("N122").Value = -9001.428400000001

$ws.Range("H136").Value = 10419506
$ws.Range("I136").Value = 12822308
$ws.Range("J136").Value = 7366.6665
$ws.Range("K136").Value = 38466924
$ws.Range("L136").Value = 22099.9995
$ws.Range("M136").Value = -38464374
$ws.Range("N136").Value = -27199.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 655.8
$ws.Range("I94").Value = 389.5
$ws.Range("J94").Value = 833.3333
$ws.Range("K94").Value = 389.5
$ws.Range("L94").Value = 833.3333
$ws.Range("M94").Value = 61.5
$ws.Range("N94").Value = -1735.3333

$ws.Range("H105").Value = 2100
$ws.Range("I105").Value = 2100
$ws.Range("K105").Value = 2100
$ws.Range("M105").Value = -353

$ws.Range("H107").Value = 202032
$ws.Range("I107").Value = 251790
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 251790
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -249870
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4451.6455
$ws.Range("I31").Value = 1628.4667
$ws.Range("J31").Value = 6180.1226
$ws.Range("K31").Value = 1628.4667
$ws.Range("L31").Value = 6180.1226
$ws.Range("M31").Value = -1333.4667
$ws.Range("N31").Value = -6770.1226

$ws.Range("H34").Value = 4451.6455
$ws.Range("I34").Value = 1628.4667
$ws.Range("J34").Value = 6180.1226
$ws.Range("K34").Value = 1628.4667
$ws.Range("L34").Value = 6180.1226
$ws.Range("M34").Value = -1426.4667
$ws.Range("N34").Value = -6584.1226

$ws.Range("H107").Value = 3677364
$ws.Range("I107").Value = 6944956.5
$ws.Range("J107").Value = 1322.5
$ws.Range("K107").Value = 6944956.5
$ws.Range("L107").Value = 1322.5
$ws.Range("M107").Value = -6943036.5
$ws.Range("N107").Value = -5162.5

$ws.Range("H122").Value = 2095.7273
$ws.Range("I122").Value = 1763.25
$ws.Range("J122").Value = 2285.7144
$ws.Range("K122").Value = 5289.75
$ws.Range("L122").Value = 6857.1432
$ws.Range("M122").Value = -2839.75
$ws.Range("N122").Value = -11757.1432

$ws.Range("H132").Value = 14707990
$ws.Range("I132").Value = 22729068
$ws.Range("J132").Value = 2682.25
$ws.Range("K132").Value = 68187204
$ws.Range("L132").Value = 8046.75
$ws.Range("M132").Value = -68184674
$ws.Range("N132").Value = -13106.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1560.2727
$ws.Range("J46").Value = 2588.3333
$ws.Range("L46").Value = 7764.999899999999
$ws.Range("N46").Value = -7946.999899999999

$ws.Range("H131").Value = 3864.3171
$ws.Range("I131").Value = 640
$ws.Range("J131").Value = 4645.9697
$ws.Range("K131").Value = 1920
$ws.Range("L131").Value = 13937.9091
$ws.Range("M131").Value = 3120
$ws.Range("N131").Value = -24017.9091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3280.8
$ws.Range("J126").Value = 3002
$ws.Range("L126").Value = 9006
$ws.Range("N126").Value = -13946

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2831.7917
$ws.Range("I61").Value = 1046.3
$ws.Range("J61").Value = 4107.143
$ws.Range("K61").Value = 1046.3
$ws.Range("L61").Value = 4107.143
$ws.Range("M61").Value = -844.3
$ws.Range("N61").Value = -4511.143

$ws.Range("H93").Value = 12598.9
$ws.Range("I93").Value = 13443.223
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 13443.223
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -12195.223
$ws.Range("N93").Value = -7496

$ws.Range("H113").Value = 2831.7917
$ws.Range("I113").Value = 1046.3
$ws.Range("J113").Value = 4107.143
$ws.Range("K113").Value = 1046.3
$ws.Range("L113").Value = 4107.143
$ws.Range("M113").Value = 1123.7
$ws.Range("N113").Value = -8447.143

$ws.Range("H122").Value = 4657.2905
$ws.Range("I122").Value = 2589.3333
$ws.Range("J122").Value = 5503.273
$ws.Range("K122").Value = 7767.999899999999
$ws.Range("L122").Value = 16509.819
$ws.Range("M122").Value = -5317.999899999999
$ws.Range("N122").Value = -21409.819

$ws.Range("H132").Value = 2557
$ws.Range("I132").Value = 2094.389
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 6283.167
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -3753.167
$ws.Range("N132").Value = -21057.9995

$ws.Range("H136").Value = 1707.8182
$ws.Range("I136").Value = 1502.1
$ws.Range("J136").Value = 3765
$ws.Range("K136").Value = 4506.299999999999
$ws.Range("L136").Value = 11295
$ws.Range("M136").Value = -1956.299999999999
$ws.Range("N136").Value = -16395

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3784.6155
$ws.Range("I96").Value = 3585.7144
$ws.Range("J96").Value = 4016.6667
$ws.Range("K96").Value = 3585.7144
$ws.Range("L96").Value = 4016.6667
$ws.Range("M96").Value = -2212.7144
$ws.Range("N96").Value = -6762.6667

$ws.Range("H122").Value = 2176.4443
$ws.Range("I122").Value = 2353.6155
$ws.Range("J122").Value = 1715.8
$ws.Range("K122").Value = 7060.8465
$ws.Range("L122").Value = 5147.4
$ws.Range("M122").Value = -4610.8465
$ws.Range("N122").Value = -10047.4

$ws.Range("H132").Value = 2978625.2
$ws.Range("I132").Value = 2133.5
$ws.Range("J132").Value = 10419854
$ws.Range("K132").Value = 6400.5
$ws.Range("L132").Value = 31259562
$ws.Range("M132").Value = -3870.5
$ws.Range("N132").Value = -31264622

$ws.Range("H136").Value = 3251.9546
$ws.Range("I136").Value = 2649
$ws.Range("J136").Value = 5302
$ws.Range("K136").Value = 7947
$ws.Range("L136").Value = 15906
$ws.Range("M136").Value = -5397
$ws.Range("N136").Value = -21006
